$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (drop the stray ".xpc" file-extension leftover from the source name)
$ws.Name = "GammaFiber2F"

# Tighten floating point precision on a couple of previously-truncated
# Gaussian-quadrature-derived intensities (row 13 and row 15)
$ws.Range("C13").Value = 0.9905057259552789
$ws.Range("D13").Value = 0.992041522875699
$ws.Range("F13").Value = 0.9905057259552789
$ws.Range("J13").Value = 0.992041522875699
$ws.Range("K13").Value = 0.991571165909284

$ws.Range("C15").Value = 0.9888161938134855
$ws.Range("F15").Value = 0.9888161938134855

# Append the new Gaussian-Quadrature "HexGrid-60degTilt5degRes" row exported
# from the new averaging scheme
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9997970327128718
$ws.Range("D16").Value = 0.9697660149721196
$ws.Range("E16").Value = 0.9997385364559017
$ws.Range("F16").Value = 0.9997970327128718
$ws.Range("G16").Value = 0.9683729531988541
$ws.Range("H16").Value = 1.001031999558077
$ws.Range("I16").Value = 0.9952941176470588
$ws.Range("J16").Value = 0.9697660149721196
$ws.Range("K16").Value = 0.9847522757140106
$ws.Range("L16").Value = 0.9922746542134413
$ws.Range("M16").Value = 0.9890001090908139

# A16 is an HKL-index cell like the other rows in column A, so it carries the
# same bold/centered/bordered style (xf index 1) -- copy it from A15 rather
# than rebuilding the format piecemeal (keeps styles.xml untouched).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
